# Auto-generated edit script: updates column F (想去人数 / "want to go" counts)
# across sheets 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types)
# to match the refreshed scrape data (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 7173
$ws.Range("F5").Value = 3
$ws.Range("F7").Value = 62
$ws.Range("F8").Value = 71
$ws.Range("F9").Value = 70
$ws.Range("F11").Value = 130
$ws.Range("F12").Value = 499
$ws.Range("F13").Value = 17
$ws.Range("F15").Value = 355
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 349
$ws.Range("F19").Value = 4068
$ws.Range("F24").Value = 1626
$ws.Range("F25").Value = 108
$ws.Range("F27").Value = 2972
$ws.Range("F28").Value = 2142
$ws.Range("F29").Value = 56
$ws.Range("F31").Value = 87
$ws.Range("F32").Value = 42
$ws.Range("F33").Value = 38
$ws.Range("F35").Value = 4195
$ws.Range("F36").Value = 440
$ws.Range("F37").Value = 315
$ws.Range("F39").Value = 942
$ws.Range("F40").Value = 755
$ws.Range("F41").Value = 166
$ws.Range("F43").Value = 1597
$ws.Range("F45").Value = 18
$ws.Range("F47").Value = 702
$ws.Range("F48").Value = 21

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 248
$ws.Range("F4").Value = 1
$ws.Range("F15").Value = 556
$ws.Range("F16").Value = 3

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 248
$ws.Range("F5").Value = 7173
$ws.Range("F9").Value = 71
$ws.Range("F10").Value = 70
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 130
$ws.Range("F14").Value = 499
$ws.Range("F17").Value = 355
$ws.Range("F18").Value = 13
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 349
$ws.Range("F21").Value = 4068
$ws.Range("F28").Value = 1626
$ws.Range("F29").Value = 108
$ws.Range("F31").Value = 2972
$ws.Range("F32").Value = 2142
$ws.Range("F33").Value = 56
$ws.Range("F35").Value = 87
$ws.Range("F36").Value = 38
$ws.Range("F39").Value = 4195
$ws.Range("F41").Value = 440
$ws.Range("F42").Value = 315
$ws.Range("F44").Value = 942
$ws.Range("F45").Value = 755
$ws.Range("F46").Value = 166
$ws.Range("F47").Value = 1597
$ws.Range("F50").Value = 702
